$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.641.97"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "2.590.64"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'307.59"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").Value = "'98.33"
$ws.Range("E6").Value = "  -2.28%  "
$ws.Range("D7").Value = "'0.593"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.573"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").Value = "'38.37"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "'53.90"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "'0.0835"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "'8.00"
$ws.Range("E13").Value = "  -3.84%  "
$ws.Range("D14").Value = "2.993.88"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "2.599.70"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "'0.904"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "'14.68"
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("D19").Value = "45.729.44"
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").Value = "'6.64"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("D22").Value = "'12.57"
$ws.Range("E22").Value = "  -5.81%  "
$ws.Range("D23").Value = "'283.51"
$ws.Range("E23").Value = "  +11.12%  "
$ws.Range("D24").Value = "'73.58"
$ws.Range("E24").Value = "  +3.48%  "
$ws.Range("D25").Value = "'3.00"
$ws.Range("E25").Value = "  -2.45%  "
$ws.Range("D26").Value = "'2.23"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "'28.92"
$ws.Range("E27").Value = "  +2.23%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'4.05"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "'10.53"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").Value = "'38.20"
$ws.Range("E31").Value = "  -6.04%  "
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("D33").Value = "'6.20"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("D35").Value = "'156.82"
$ws.Range("E35").Value = "  +2.45%  "
$ws.Range("D36").Value = "'2.23"
$ws.Range("E36").Value = "  -2.76%  "
$ws.Range("E37").Value = "  -3.04%  "
$ws.Range("D38").Value = "'0.0825"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").Value = "'0.122"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "'15.83"
$ws.Range("E41").Value = "  -7.26%  "
$ws.Range("D42").Value = "'0.0323"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").Value = "'21.36"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("D44").Value = "'3.50"
$ws.Range("E44").Value = "  -3.57%  "
$ws.Range("D45").Value = "'3.97"
$ws.Range("E45").Value = "  -6.59%  "
$ws.Range("D46").Value = "2.098.82"
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("D47").Value = "'0.999"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "'93.51"
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("D49").Value = "'9.17"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").Value = "'107.92"
$ws.Range("E50").Value = "  -3.04%  "
$ws.Range("D51").Value = "2.846.58"
$ws.Range("E51").Value = "  -0.78%  "
